$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "수학을 잘하는 법 (SIAI기준)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/be-good-at-math-by-siai/#utm_source=rss&utm_medium=rss&utm_campaign=be-good-at-math-by-siai"

$ws.Range("D37").Value = "[Paper Review] Simple Unsupervised Keyphrase Extraction using Sentence Embeddings (EmbedRank)"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1947&mod=document&pageid=1"
